$wb = $excel.ActiveWorkbook

# The "Italy" sheet is the template for the new "Spain" sheet.
$italy = $wb.Worksheets.Item("Italy")

# Duplicate Italy to create the new sheet right after it. The copy
# inherits Italy's current selection state (cell B4 was selected).
$italy.Copy($null, $italy)

# The copy becomes the active sheet, positioned immediately after "Italy".
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Update the market name and product code cells for Spain.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2064/T2063"

# Select the whole data range on Italy (its tab is no longer the active
# one), then re-activate Spain so it remains the active/selected tab.
$italy.Activate()
$italy.Range("A1:D13").Select()
$spain.Activate()
